$d = $word.ActiveDocument

# --- Delete the paragraph containing "See post reactions." ---
$found = $false
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*See post reactions.*") {
        $p.Range.Delete()
        $found = $true
    }
}
Write-Output ("Deleted 'See post reactions.' paragraph: " + $found)

# --- Delete the empty paragraph right after "Search songs functionality." ---
$found2 = $false
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Search songs functionality.*") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text.Trim() -eq "") {
            $next.Range.Delete()
            $found2 = $true
        }
    }
}
Write-Output ("Deleted empty paragraph after 'Search songs functionality.': " + $found2)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
